# Update the "想去人数" (want-to-go count) column F values on both the
# "展览" and "全部类型" worksheets, which contain duplicated data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 151
    3  = 1704
    4  = 787
    7  = 11920
    8  = 39
    10 = 475
    12 = 1111
    13 = 845
    14 = 13459
    15 = 13421
    17 = 153
    20 = 281
    21 = 95
    23 = 71
    24 = 165
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
